$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transaction Service")

# Row 2 title changes to "Transaction Service"
$ws.Range("A2").Value = "Transaction Service"

# Row 4
$ws.Range("A4").Value = "R001"
$ws.Range("B4").Value = "T001"
$ws.Range("C4").Value = "saveTransaction(Transaction transaction)"
$ws.Range("D4").Value = $true
$ws.Range("F4").Value = "Add transaction "

# Row 7
$ws.Range("A7").Value = "R002"
$ws.Range("B7").Value = "T001"
$ws.Range("C7").Value = "saveTransactionByUserName(String name)"
$ws.Range("D7").Value = $true
$ws.Range("F7").Value = "if any name is passed"

# Row 8
$ws.Range("A8").Value = "R002"
$ws.Range("B8").Value = "T002"
$ws.Range("C8").Value = "saveTransactionByUserName(null)"
$ws.Range("D8").Value = $false
$ws.Range("F8").Value = "null object is passed"

# Row 10
$ws.Range("A10").Value = "R003"
$ws.Range("B10").Value = "T001"
$ws.Range("C10").Value = "getMaxTransactionId()"
$ws.Range("D10").Value = "maximum no.of transactions done by customer"

# Row 12
$ws.Range("A12").Value = "R004"
$ws.Range("B12").Value = "T001"
$ws.Range("C12").Value = "getAllTransaction()"
$ws.Range("D12").Value = "List<Transactions>"
$ws.Range("F12").Value = "if atleast one transaction exists"

# Row 13
$ws.Range("A13").Value = "R004"
$ws.Range("B13").Value = "T002"
$ws.Range("C13").Value = "getAllTransaction()"
$ws.Range("D13").Value = "null"
$ws.Range("F13").Value = "if no transaction is done"

# Row 15 / 16 (string pool order requires both C cells before the F cells)
$ws.Range("A15").Value = "R005"
$ws.Range("B15").Value = "T001"
$ws.Range("C15").Value = "getAllTransactionByName(String name)"
$ws.Range("D15").Value = "List<Transactions>"

$ws.Range("A16").Value = "R005"
$ws.Range("B16").Value = "T002"
$ws.Range("C16").Value = "getAllTransaction(String name)"
$ws.Range("D16").Value = "null"

$ws.Range("F15").Value = "if atleast one transaction exists for customer"
$ws.Range("F16").Value = "if no transaction is done for customer"

# Column widths (best-fit, matching the auto-sized columns from the authored file)
$ws.Columns.Item(1).ColumnWidth = 14
$ws.Columns.Item(2).ColumnWidth = 10.166666666666666
$ws.Columns.Item(3).ColumnWidth = 34.333333333333336
$ws.Columns.Item(4).ColumnWidth = 39.333333333333336
$ws.Columns.Item(5).ColumnWidth = 12.5
$ws.Columns.Item(6).ColumnWidth = 37.333333333333336

# View settings
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("F16").Select()
